$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.319.07"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "3.380.30"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "576.36"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "134.22"
$ws.Range("E6").Value = "  +5.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.381.05"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "7.60"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "3.953.04"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "3.377.88"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "25.20"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "61.417.82"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "13.98"
$ws.Range("E19").Value = "  +5.81%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "9.36"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "379.51"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").Value = "3.511.64"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "70.74"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +9.20%  "
$ws.Range("D28").Value = "1.68"
$ws.Range("E28").Value = "  +18.68%  "
$ws.Range("D29").Value = "7.73"
$ws.Range("E29").Value = "  +10.12%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "0.156"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "3.412.87"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "23.41"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").Value = "5.58"
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "1.56"
$ws.Range("D40").Value = "162.60"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "1.23"
$ws.Range("E43").Value = "  +12.10%  "
$ws.Range("D44").Value = "4.43"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").Value = "41.68"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "0.758"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("D48").Value = "23.67"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").Value = "23.18"
$ws.Range("E50").Value = "  +13.44%  "
$ws.Range("D51").Value = "0.902"
$ws.Range("E51").Value = "  +4.02%  "
